$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.564.12"
$ws.Range("E2").Value = "  +2.54%  "

# Row 3
$ws.Range("D3").Value = "3.813.08"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "686.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.15%  "

# Row 7
$ws.Range("D7").Value = "3.812.42"
$ws.Range("E7").Value = "  +1.27%  "

# Row 9
$ws.Range("E9").Value = "  +1.09%  "

# Row 10
$ws.Range("E10").Value = "  +2.06%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.09%  "

# Row 13
$ws.Range("E13").Value = "  +0.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "

# Row 15
$ws.Range("D15").Value = "4.454.90"
$ws.Range("E15").Value = "  +1.24%  "

# Row 16
$ws.Range("D16").Value = "3.810.54"
$ws.Range("E16").Value = "  +2.09%  "

# Row 17
$ws.Range("D17").Value = "70.639.88"
$ws.Range("E17").Value = "  +2.69%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.114"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.62%  "

# Row 23
$ws.Range("E23").Value = "  +1.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("E25").Value = "  -0.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("D30").Value = "3.963.20"
$ws.Range("E30").Value = "  +1.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.37%  "

# Row 34
$ws.Range("E34").Value = "  +3.81%  "

# Row 35
$ws.Range("E35").Value = "  +6.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.73%  "

# Row 37
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").Value = "3.762.59"
$ws.Range("E38").Value = "  +1.23%  "

# Row 39
$ws.Range("E39").Value = "  +1.48%  "

# Row 40
$ws.Range("E40").Value = "  +4.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.32%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.95%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "

# Row 45
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("E46").Value = "  +7.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.05%  "

# Row 51
$ws.Range("E51").Value = "  +2.02%  "
